$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename two existing rows (correcting copy/typo'd labels)
$ws.Range("C9").Value = "Windows Events"
$ws.Range("C10").Value = "Schedule Task"

# Add two new IOC control rows at the bottom of the table. Seed the new
# rows by copying the formatting of the row above (keeps the same cell
# styles as the rest of the table) and then overwrite the values.
$ws.Range("A15:C15").Copy($ws.Range("A16:C16"))
$ws.Range("A15:C15").Copy($ws.Range("A17:C17"))

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Windows"

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "Windows"

# (Configuration Files is entered before List All User Accounts so the
# shared-string table ends up in the same add-order as the source file.)
$ws.Range("C17").Value = "Configuration Files"
$ws.Range("C16").Value = "List All User Accounts"

# Update selection to reflect last-edited cell
$ws.Range("C17").Select()
